$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix C4 total price (was 36, should be 0)
$ws.Range("C4").Value = 0

# New row 8: another order from Elmar Qarayev, now Canceled
$ws.Range("A8").Value = "Elmar Qarayev"
$ws.Range("B8").Value = "elmarqarayev69@gmail.com"
$ws.Range("C8").Value = 396
$ws.Range("D8").Value = "Canceled"

# New row 9: order from new customer Elmar Garayev, Pending
$ws.Range("A9").Value = "Elmar Garayev"
$ws.Range("B9").Value = "elmareg@code.edu.az"
$ws.Range("C9").Value = 156
$ws.Range("D9").Value = "Pending"

# New row 10: another order from Elmar Garayev, Pending
$ws.Range("A10").Value = "Elmar Garayev"
$ws.Range("B10").Value = "elmareg@code.edu.az"
$ws.Range("C10").Value = 66
$ws.Range("D10").Value = "Pending"
